# Update column G ("K") values on the active sheet to reflect the
# regenerated strikeout counts (replacing the old "Strike#" values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 3
    17 = 2
    18 = 2
    19 = 2
    20 = 3
    21 = 0
    22 = 0
    23 = 4
    24 = 1
    25 = 0
    26 = 2
    27 = 0
    28 = 2
    29 = 1
    30 = 5
    31 = 4
    32 = 2
    33 = 1
    34 = 2
    35 = 2
    36 = 3
    37 = 0
    38 = 3
    39 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
